$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# There are two distinct logos reused across the "default" and "first page"
# headers/footers of this section:
#   - the BTEC logo (header1.xml / header2.xml) currently named "image2.jpg"
#     needs to become "image1.jpg"
#   - the Pearson logo (footer1.xml / footer2.xml) currently named "image1.png"
#     needs to become "image2.png"
# (the picture's alternative text / description is untouched - only the
# shape's name changes)

function Rename-InlineShape($range, $newName) {
    $ishp = $range.InlineShapes.Item(1)
    $shp = $ishp.ConvertToShape()
    $shp.Name = $newName
    $shp.ConvertToInlineShape() | Out-Null
}

# Headers: wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
Rename-InlineShape $sec.Headers.Item(1).Range "image1.jpg"
Rename-InlineShape $sec.Headers.Item(2).Range "image1.jpg"

# Footers: wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
Rename-InlineShape $sec.Footers.Item(1).Range "image2.png"
Rename-InlineShape $sec.Footers.Item(2).Range "image2.png"
